$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date (column C) for rows 2-5 from 2023-10-13 (45212) to 2023-10-22 (45221)
$newDate = Get-Date -Year 2023 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2:C5").Value = $newDate
